$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "masthoogte" in J1 (new attribute/column)
$ws.Range("J1").Value = "masthoogte"

# Resize the columns that now have explicit widths in the target sheet (G, H, I, J)
$ws.Columns.Item(7).ColumnWidth = 9.333333333333334
$ws.Columns.Item(8).ColumnWidth = 15.666666666666666
$ws.Columns.Item(9).ColumnWidth = 22.333333333333332
$ws.Columns.Item(10).ColumnWidth = 10.666666666666666

# Move/resize the current selection to I2:I6 with I2 as the active cell
$ws.Range("I2:I6").Select()

$wb.Save()
